$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Sparse odds updates for existing rows ----
# Row 2
$ws.Cells.Item(2, 8).Value = 4.2
$ws.Cells.Item(2, 9).Value = 6.25
$ws.Cells.Item(2, 14).Value = 12
$ws.Cells.Item(2, 21).Value = 1.8
$ws.Cells.Item(2, 22).Value = 1.95
$ws.Cells.Item(2, 29).Value = 12
$ws.Cells.Item(2, 52).Value = 101

# Row 3
$ws.Cells.Item(3, 7).Value = 2.37
$ws.Cells.Item(3, 8).Value = 3.6
$ws.Cells.Item(3, 9).Value = 2.62
$ws.Cells.Item(3, 10).Value = 2.87
$ws.Cells.Item(3, 11).Value = 2.27
$ws.Cells.Item(3, 12).Value = 3.1
$ws.Cells.Item(3, 15).Value = 1.2
$ws.Cells.Item(3, 16).Value = 4
$ws.Cells.Item(3, 17).Value = 1.62
$ws.Cells.Item(3, 18).Value = 2.2
$ws.Cells.Item(3, 19).Value = 1.31
$ws.Cells.Item(3, 20).Value = 3.15
$ws.Cells.Item(3, 21).Value = 1.53
$ws.Cells.Item(3, 23).Value = 10.75
$ws.Cells.Item(3, 24).Value = 14
$ws.Cells.Item(3, 25).Value = 9.25
$ws.Cells.Item(3, 26).Value = 26
$ws.Cells.Item(3, 27).Value = 17
$ws.Cells.Item(3, 28).Value = 22
$ws.Cells.Item(3, 30).Value = 7.2
$ws.Cells.Item(3, 35).Value = 15
$ws.Cells.Item(3, 36).Value = 9.75
$ws.Cells.Item(3, 37).Value = 30
$ws.Cells.Item(3, 38).Value = 19.5
$ws.Cells.Item(3, 39).Value = 24
$ws.Cells.Item(3, 40).Value = 4.6
$ws.Cells.Item(3, 41).Value = 12
$ws.Cells.Item(3, 42).Value = 17
$ws.Cells.Item(3, 43).Value = 45
$ws.Cells.Item(3, 44).Value = 65
$ws.Cells.Item(3, 46).Value = 3.15
$ws.Cells.Item(3, 49).Value = 4.8
$ws.Cells.Item(3, 50).Value = 13.5
$ws.Cells.Item(3, 51).Value = 18.5
$ws.Cells.Item(3, 52).Value = 55
$ws.Cells.Item(3, 53).Value = 75
$ws.Cells.Item(3, 54).Value = 175

# Row 4
$ws.Cells.Item(4, 7).Value = 2.12
$ws.Cells.Item(4, 8).Value = 3.15
$ws.Cells.Item(4, 9).Value = 3.45
$ws.Cells.Item(4, 10).Value = 2.7
$ws.Cells.Item(4, 11).Value = 2.07
$ws.Cells.Item(4, 12).Value = 3.9
$ws.Cells.Item(4, 14).Value = 6.5
$ws.Cells.Item(4, 20).Value = 2.65
$ws.Cells.Item(4, 21).Value = 1.85
$ws.Cells.Item(4, 22).Value = 1.87
$ws.Cells.Item(4, 23).Value = 6.8
$ws.Cells.Item(4, 24).Value = 9.75
$ws.Cells.Item(4, 25).Value = 8.75
$ws.Cells.Item(4, 26).Value = 20
$ws.Cells.Item(4, 27).Value = 18
$ws.Cells.Item(4, 29).Value = 6.5
$ws.Cells.Item(4, 30).Value = 6
$ws.Cells.Item(4, 31).Value = 15
$ws.Cells.Item(4, 32).Value = 75
$ws.Cells.Item(4, 33).Value = 700
$ws.Cells.Item(4, 34).Value = 9
$ws.Cells.Item(4, 35).Value = 17.5
$ws.Cells.Item(4, 36).Value = 11.75
$ws.Cells.Item(4, 37).Value = 50
$ws.Cells.Item(4, 38).Value = 32
$ws.Cells.Item(4, 40).Value = 4
$ws.Cells.Item(4, 41).Value = 11
$ws.Cells.Item(4, 42).Value = 19
$ws.Cells.Item(4, 44).Value = 75
$ws.Cells.Item(4, 46).Value = 2.65
$ws.Cells.Item(4, 47).Value = 7
$ws.Cells.Item(4, 48).Value = 65
$ws.Cells.Item(4, 49).Value = 5.3
$ws.Cells.Item(4, 50).Value = 19
$ws.Cells.Item(4, 51).Value = 26
$ws.Cells.Item(4, 52).Value = 100
$ws.Cells.Item(4, 53).Value = 150
$ws.Cells.Item(4, 54).Value = 350

# Row 5
$ws.Cells.Item(5, 7).Value = 2.65
$ws.Cells.Item(5, 8).Value = 4
$ws.Cells.Item(5, 9).Value = 2.22
$ws.Cells.Item(5, 10).Value = 3
$ws.Cells.Item(5, 11).Value = 2.52
$ws.Cells.Item(5, 12).Value = 2.65
$ws.Cells.Item(5, 17).Value = 1.38
$ws.Cells.Item(5, 18).Value = 2.8
$ws.Cells.Item(5, 21).Value = 1.38
$ws.Cells.Item(5, 22).Value = 2.82
$ws.Cells.Item(5, 23).Value = 16.5
$ws.Cells.Item(5, 24).Value = 19.5
$ws.Cells.Item(5, 25).Value = 10.75
$ws.Cells.Item(5, 26).Value = 32
$ws.Cells.Item(5, 27).Value = 18
$ws.Cells.Item(5, 28).Value = 19
$ws.Cells.Item(5, 29).Value = 10
$ws.Cells.Item(5, 31).Value = 11
$ws.Cells.Item(5, 32).Value = 29
$ws.Cells.Item(5, 36).Value = 9.75
$ws.Cells.Item(5, 37).Value = 25
$ws.Cells.Item(5, 38).Value = 15
$ws.Cells.Item(5, 39).Value = 17.5
$ws.Cells.Item(5, 40).Value = 5.3
$ws.Cells.Item(5, 41).Value = 13
$ws.Cells.Item(5, 42).Value = 15
$ws.Cells.Item(5, 43).Value = 45
$ws.Cells.Item(5, 44).Value = 55
$ws.Cells.Item(5, 47).Value = 6.2
$ws.Cells.Item(5, 49).Value = 4.85
$ws.Cells.Item(5, 50).Value = 10.75
$ws.Cells.Item(5, 51).Value = 14
$ws.Cells.Item(5, 52).Value = 35
$ws.Cells.Item(5, 53).Value = 45
$ws.Cells.Item(5, 54).Value = 110

# Row 6
$ws.Cells.Item(6, 7).Value = 2.87
$ws.Cells.Item(6, 8).Value = 3.5
$ws.Cells.Item(6, 9).Value = 2.25
$ws.Cells.Item(6, 11).Value = 2.22
$ws.Cells.Item(6, 12).Value = 2.8
$ws.Cells.Item(6, 14).Value = 8.25
$ws.Cells.Item(6, 16).Value = 3.8
$ws.Cells.Item(6, 17).Value = 1.7
$ws.Cells.Item(6, 18).Value = 2.07
$ws.Cells.Item(6, 21).Value = 1.57
$ws.Cells.Item(6, 22).Value = 2.25
$ws.Cells.Item(6, 23).Value = 11.5
$ws.Cells.Item(6, 24).Value = 17
$ws.Cells.Item(6, 26).Value = 35
$ws.Cells.Item(6, 27).Value = 22
$ws.Cells.Item(6, 28).Value = 25
$ws.Cells.Item(6, 29).Value = 8.25
$ws.Cells.Item(6, 30).Value = 6.9
$ws.Cells.Item(6, 31).Value = 12
$ws.Cells.Item(6, 34).Value = 9.5
$ws.Cells.Item(6, 37).Value = 23
$ws.Cells.Item(6, 38).Value = 16.5
$ws.Cells.Item(6, 40).Value = 5
$ws.Cells.Item(6, 41).Value = 15
$ws.Cells.Item(6, 42).Value = 19.5
$ws.Cells.Item(6, 45).Value = 200
$ws.Cells.Item(6, 47).Value = 6.6
$ws.Cells.Item(6, 49).Value = 4.35
$ws.Cells.Item(6, 50).Value = 11.5
$ws.Cells.Item(6, 52).Value = 45
$ws.Cells.Item(6, 53).Value = 70

# Row 8
$ws.Cells.Item(8, 7).Value = 2.1
$ws.Cells.Item(8, 8).Value = 3.25
$ws.Cells.Item(8, 9).Value = 3.6
$ws.Cells.Item(8, 10).Value = 2.75
$ws.Cells.Item(8, 14).Value = 9.5
$ws.Cells.Item(8, 26).Value = 19
$ws.Cells.Item(8, 31).Value = 15
$ws.Cells.Item(8, 36).Value = 13

# Row 10
$ws.Cells.Item(10, 7).Value = 1.57
$ws.Cells.Item(10, 9).Value = 5.5
$ws.Cells.Item(10, 11).Value = 2.3
$ws.Cells.Item(10, 17).Value = 1.8
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 26).Value = 12
$ws.Cells.Item(10, 27).Value = 13
$ws.Cells.Item(10, 31).Value = 17
$ws.Cells.Item(10, 33).Value = 251
$ws.Cells.Item(10, 34).Value = 15
$ws.Cells.Item(10, 36).Value = 17
$ws.Cells.Item(10, 37).Value = 51
$ws.Cells.Item(10, 41).Value = 8

# Row 12
$ws.Cells.Item(12, 7).Value = 1.9
$ws.Cells.Item(12, 9).Value = 3.8
$ws.Cells.Item(12, 13).Value = 1.04
$ws.Cells.Item(12, 14).Value = 13
$ws.Cells.Item(12, 17).Value = 1.7
$ws.Cells.Item(12, 18).Value = 2.1
$ws.Cells.Item(12, 24).Value = 10
$ws.Cells.Item(12, 27).Value = 13
$ws.Cells.Item(12, 38).Value = 29
$ws.Cells.Item(12, 40).Value = 4
$ws.Cells.Item(12, 43).Value = 29
$ws.Cells.Item(12, 54).Value = 151

# Row 13
$ws.Cells.Item(13, 17).Value = 1.67
$ws.Cells.Item(13, 18).Value = 2.15

# Row 16
$ws.Cells.Item(16, 7).Value = 1.67
$ws.Cells.Item(16, 8).Value = 3.4
$ws.Cells.Item(16, 9).Value = 4.75
$ws.Cells.Item(16, 10).Value = 2.4
$ws.Cells.Item(16, 11).Value = 2.05
$ws.Cells.Item(16, 12).Value = 6
$ws.Cells.Item(16, 19).Value = 1.5
$ws.Cells.Item(16, 20).Value = 2.5
$ws.Cells.Item(16, 28).Value = 34
$ws.Cells.Item(16, 46).Value = 2.5

# Row 17
$ws.Cells.Item(17, 7).Value = 2.2
$ws.Cells.Item(17, 9).Value = 3
$ws.Cells.Item(17, 10).Value = 3.25
$ws.Cells.Item(17, 12).Value = 4.33
$ws.Cells.Item(17, 13).Value = 1.1
$ws.Cells.Item(17, 14).Value = 7
$ws.Cells.Item(17, 15).Value = 1.53
$ws.Cells.Item(17, 16).Value = 2.38
$ws.Cells.Item(17, 23).Value = 5.5
$ws.Cells.Item(17, 24).Value = 9.5
$ws.Cells.Item(17, 25).Value = 10
$ws.Cells.Item(17, 26).Value = 21
$ws.Cells.Item(17, 27).Value = 23
$ws.Cells.Item(17, 34).Value = 7
$ws.Cells.Item(17, 35).Value = 15
$ws.Cells.Item(17, 36).Value = 13
$ws.Cells.Item(17, 38).Value = 34
$ws.Cells.Item(17, 39).Value = 51
$ws.Cells.Item(17, 40).Value = 4
$ws.Cells.Item(17, 42).Value = 29
$ws.Cells.Item(17, 44).Value = 81
$ws.Cells.Item(17, 45).Value = 301
$ws.Cells.Item(17, 49).Value = 5
$ws.Cells.Item(17, 50).Value = 21
$ws.Cells.Item(17, 51).Value = 41
$ws.Cells.Item(17, 52).Value = 81

# Row 18
$ws.Cells.Item(18, 7).Value = 2.2
$ws.Cells.Item(18, 9).Value = 3.2
$ws.Cells.Item(18, 12).Value = 3.6
$ws.Cells.Item(18, 36).Value = 12

# Row 23
$ws.Cells.Item(23, 7).Value = 4.15
$ws.Cells.Item(23, 9).Value = 1.78
$ws.Cells.Item(23, 10).Value = 4.45
$ws.Cells.Item(23, 12).Value = 2.35
$ws.Cells.Item(23, 23).Value = 11.5
$ws.Cells.Item(23, 24).Value = 25
$ws.Cells.Item(23, 25).Value = 14.5
$ws.Cells.Item(23, 26).Value = 70
$ws.Cells.Item(23, 27).Value = 40
$ws.Cells.Item(23, 28).Value = 45
$ws.Cells.Item(23, 34).Value = 7.1
$ws.Cells.Item(23, 35).Value = 9
$ws.Cells.Item(23, 37).Value = 15.5
$ws.Cells.Item(23, 38).Value = 15
$ws.Cells.Item(23, 40).Value = 5.9
$ws.Cells.Item(23, 41).Value = 23
$ws.Cells.Item(23, 43).Value = 120
$ws.Cells.Item(23, 46).Value = 2.82
$ws.Cells.Item(23, 49).Value = 3.65
$ws.Cells.Item(23, 50).Value = 8.75
$ws.Cells.Item(23, 52).Value = 30
$ws.Cells.Item(23, 53).Value = 65

# ---- Full row replacements (new/moved matches) ----
# Row 20
$row20 = New-Object 'object[,]' 1,56
$row20[0,0] = "fyqcIbbe"
$row20[0,1] = "22/11/2024"
$row20[0,2] = "15:30"
$row20[0,3] = "ROMANIA - LIGA 1"
$row20[0,4] = "Farul Constanta"
$row20[0,5] = "Otelul"
$row20[0,6] = 2.35
$row20[0,7] = 2.9
$row20[0,8] = 3
$row20[0,9] = 3.25
$row20[0,10] = 1.91
$row20[0,11] = 4
$row20[0,12] = 1.1
$row20[0,13] = 7
$row20[0,14] = 1.5
$row20[0,15] = 2.5
$row20[0,16] = 2.5
$row20[0,17] = 1.5
$row20[0,18] = 1.57
$row20[0,19] = 2.25
$row20[0,20] = 2.1
$row20[0,21] = 1.67
$row20[0,22] = 6.5
$row20[0,23] = 10
$row20[0,24] = 10
$row20[0,25] = 23
$row20[0,26] = 23
$row20[0,27] = 41
$row20[0,28] = 6.5
$row20[0,29] = 6
$row20[0,30] = 19
$row20[0,31] = 67
$row20[0,32] = 1000
$row20[0,33] = 7.5
$row20[0,34] = 13
$row20[0,35] = 12
$row20[0,36] = 34
$row20[0,37] = 29
$row20[0,38] = 41
$row20[0,39] = 4.33
$row20[0,40] = 15
$row20[0,41] = 29
$row20[0,42] = 51
$row20[0,43] = 81
$row20[0,44] = 251
$row20[0,45] = 2.25
$row20[0,46] = 9
$row20[0,47] = 67
$row20[0,48] = 5
$row20[0,49] = 19
$row20[0,50] = 34
$row20[0,51] = 67
$row20[0,52] = 101
$row20[0,53] = 301
$row20[0,54] = 51
$row20[0,55] = 51
$ws.Range("A20:BD20").Value = $row20

# Row 21
$row21 = New-Object 'object[,]' 1,56
$row21[0,0] = "xYNxflsK"
$row21[0,1] = "22/11/2024"
$row21[0,2] = "14:00"
$row21[0,3] = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$row21[0,4] = "Al Nassr"
$row21[0,5] = "Al Qadisiya"
$row21[0,6] = 1.65
$row21[0,7] = 4
$row21[0,8] = 4.2
$row21[0,9] = 2.2
$row21[0,10] = 2.3
$row21[0,11] = 4.5
$row21[0,12] = 1.03
$row21[0,13] = 10
$row21[0,14] = 1.2
$row21[0,15] = 4.33
$row21[0,16] = 1.67
$row21[0,17] = 2.15
$row21[0,18] = 1.33
$row21[0,19] = 3.25
$row21[0,20] = 1.73
$row21[0,21] = 2
$row21[0,22] = 8.5
$row21[0,23] = 9
$row21[0,24] = 8.5
$row21[0,25] = 13
$row21[0,26] = 13
$row21[0,27] = 23
$row21[0,28] = 15
$row21[0,29] = 8
$row21[0,30] = 15
$row21[0,31] = 41
$row21[0,32] = 151
$row21[0,33] = 15
$row21[0,34] = 23
$row21[0,35] = 15
$row21[0,36] = 51
$row21[0,37] = 34
$row21[0,38] = 34
$row21[0,39] = 3.75
$row21[0,40] = 8.5
$row21[0,41] = 17
$row21[0,42] = 23
$row21[0,43] = 41
$row21[0,44] = 101
$row21[0,45] = 3.25
$row21[0,46] = 8
$row21[0,47] = 51
$row21[0,48] = 6.5
$row21[0,49] = 23
$row21[0,50] = 29
$row21[0,51] = 81
$row21[0,52] = 81
$row21[0,53] = 151
$row21[0,54] = 81
$row21[0,55] = 81
$ws.Range("A21:BD21").Value = $row21

# Row 25
$row25 = New-Object 'object[,]' 1,56
$row25[0,0] = "WWp90WXG"
$row25[0,1] = "22/11/2024"
$row25[0,2] = "15:30"
$row25[0,3] = "SWITZERLAND - CHALLENGE LEAGUE"
$row25[0,4] = "Schaffhausen"
$row25[0,5] = "Lausanne Ouchy"
$row25[0,6] = 3.1
$row25[0,7] = 3.25
$row25[0,8] = 2.22
$row25[0,9] = 3.6
$row25[0,10] = 2.05
$row25[0,11] = 2.9
$row25[0,12] = 1.06
$row25[0,13] = 7.5
$row25[0,14] = 1.27
$row25[0,15] = 3.4
$row25[0,16] = 1.82
$row25[0,17] = 1.93
$row25[0,18] = 1.42
$row25[0,19] = 2.65
$row25[0,20] = 1.65
$row25[0,21] = 2.12
$row25[0,22] = 11
$row25[0,23] = 17.5
$row25[0,24] = 10.5
$row25[0,25] = 40
$row25[0,26] = 24
$row25[0,27] = 28
$row25[0,28] = 7.5
$row25[0,29] = 6.3
$row25[0,30] = 12.5
$row25[0,31] = 50
$row25[0,32] = 350
$row25[0,33] = 8.25
$row25[0,34] = 11.25
$row25[0,35] = 8.75
$row25[0,36] = 22
$row25[0,37] = 17.5
$row25[0,38] = 25
$row25[0,39] = 5
$row25[0,40] = 17
$row25[0,41] = 23
$row25[0,42] = 80
$row25[0,43] = 110
$row25[0,44] = 300
$row25[0,45] = 2.65
$row25[0,46] = 6.9
$row25[0,47] = 60
$row25[0,48] = 4.2
$row25[0,49] = 12.5
$row25[0,50] = 21
$row25[0,51] = 50
$row25[0,52] = 90
$row25[0,53] = 250
$row25[0,54] = 51
$row25[0,55] = 51
$ws.Range("A25:BD25").Value = $row25

# Row 26
$row26 = New-Object 'object[,]' 1,56
$row26[0,0] = "4tmXLw0f"
$row26[0,1] = "22/11/2024"
$row26[0,2] = "15:30"
$row26[0,3] = "SWITZERLAND - CHALLENGE LEAGUE"
$row26[0,4] = "Stade Nyonnais"
$row26[0,5] = "Vaduz"
$row26[0,6] = 2.67
$row26[0,7] = 3.55
$row26[0,8] = 2.35
$row26[0,9] = 3.15
$row26[0,10] = 2.3
$row26[0,11] = 2.82
$row26[0,12] = 1.03
$row26[0,13] = 9.25
$row26[0,14] = 1.16
$row26[0,15] = 4.6
$row26[0,16] = 1.5
$row26[0,17] = 2.42
$row26[0,18] = 1.29
$row26[0,19] = 3.3
$row26[0,20] = 1.42
$row26[0,21] = 2.65
$row26[0,22] = 13
$row26[0,23] = 17.5
$row26[0,24] = 10
$row26[0,25] = 32
$row26[0,26] = 19
$row26[0,27] = 21
$row26[0,28] = 9.25
$row26[0,29] = 7.5
$row26[0,30] = 10.75
$row26[0,31] = 32
$row26[0,32] = 150
$row26[0,33] = 13
$row26[0,34] = 15.5
$row26[0,35] = 9.25
$row26[0,36] = 27
$row26[0,37] = 16.5
$row26[0,38] = 18.5
$row26[0,39] = 5.1
$row26[0,40] = 14
$row26[0,41] = 17
$row26[0,42] = 55
$row26[0,43] = 70
$row26[0,44] = 150
$row26[0,45] = 3.3
$row26[0,46] = 6.1
$row26[0,47] = 37
$row26[0,48] = 4.75
$row26[0,49] = 11.75
$row26[0,50] = 15.5
$row26[0,51] = 40
$row26[0,52] = 55
$row26[0,53] = 150
$row26[0,54] = 500
$row26[0,55] = 51
$ws.Range("A26:BD26").Value = $row26

# Row 27
$row27 = New-Object 'object[,]' 1,56
$row27[0,0] = "jPRoB7i2"
$row27[0,1] = "22/11/2024"
$row27[0,2] = "14:00"
$row27[0,3] = "TURKEY - 1. LIG"
$row27[0,4] = "Pendikspor"
$row27[0,5] = "Umraniyespor"
$row27[0,6] = 1.85
$row27[0,7] = 3.6
$row27[0,8] = 3.9
$row27[0,9] = 2.5
$row27[0,10] = 2.1
$row27[0,11] = 4.75
$row27[0,12] = 1.06
$row27[0,13] = 10
$row27[0,14] = 1.33
$row27[0,15] = 3.25
$row27[0,16] = 2.05
$row27[0,17] = 1.75
$row27[0,18] = 1.44
$row27[0,19] = 2.63
$row27[0,20] = 1.83
$row27[0,21] = 1.83
$row27[0,22] = 7
$row27[0,23] = 8.5
$row27[0,24] = 9
$row27[0,25] = 15
$row27[0,26] = 15
$row27[0,27] = 29
$row27[0,28] = 9.5
$row27[0,29] = 7
$row27[0,30] = 17
$row27[0,31] = 51
$row27[0,32] = 301
$row27[0,33] = 10
$row27[0,34] = 21
$row27[0,35] = 13
$row27[0,36] = 41
$row27[0,37] = 34
$row27[0,38] = 41
$row27[0,39] = 3.75
$row27[0,40] = 10
$row27[0,41] = 21
$row27[0,42] = 34
$row27[0,43] = 51
$row27[0,44] = 151
$row27[0,45] = 2.63
$row27[0,46] = 8.5
$row27[0,47] = 51
$row27[0,48] = 6
$row27[0,49] = 23
$row27[0,50] = 34
$row27[0,51] = 81
$row27[0,52] = 101
$row27[0,53] = 251
$row27[0,54] = 126
$row27[0,55] = 126
$ws.Range("A27:BD27").Value = $row27

